$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2483.3333
$ws.Range("I98").Value = 2386.5386
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 2386.5386
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -888.5385999999999
$ws.Range("N98").Value = -7996
$ws.Range("H122").Value = 2483.3333
$ws.Range("I122").Value = 2386.5386
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7159.6158
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4709.6158
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 259376
$ws.Range("I132").Value = 259376
$ws.Range("K132").Value = 778128
$ws.Range("M132").Value = -775598
$ws.Range("H137").Value = 4703.243
$ws.Range("I137").Value = 7063.0557
$ws.Range("J137").Value = 2467.6316
$ws.Range("K137").Value = 21189.1671
$ws.Range("L137").Value = 7402.8948
$ws.Range("M137").Value = -18639.1671
$ws.Range("N137").Value = -12502.8948
$ws.Range("H141").Value = 2772.6191
$ws.Range("I141").Value = 2772.6191
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8317.8573
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3137.8573
$ws.Range("N141").Value = $null
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2674228.2
$ws.Range("I2").Value = 491.5
$ws.Range("J2").Value = 5882712.5
$ws.Range("K2").Value = 491.5
$ws.Range("L2").Value = 5882712.5
$ws.Range("M2").Value = -378.5
$ws.Range("N2").Value = -5882938.5
$ws.Range("H32").Value = 6767.1025
$ws.Range("I32").Value = 3753.9194
$ws.Range("K32").Value = 3753.9194
$ws.Range("M32").Value = -3466.9194
$ws.Range("H61").Value = 1913.4
$ws.Range("I61").Value = 1682.8889
$ws.Range("J61").Value = 2506.1428
$ws.Range("K61").Value = 1682.8889
$ws.Range("L61").Value = 2506.1428
$ws.Range("M61").Value = -1470.8889
$ws.Range("N61").Value = -2930.1428
$ws.Range("H74").Value = 2781468
$ws.Range("I74").Value = 4348425
$ws.Range("J74").Value = 9159.462
$ws.Range("K74").Value = 4348425
$ws.Range("L74").Value = 9159.462
$ws.Range("M74").Value = -4347551
$ws.Range("N74").Value = -10907.462
$ws.Range("H77").Value = 2781468
$ws.Range("I77").Value = 4348425
$ws.Range("J77").Value = 9159.462
$ws.Range("K77").Value = 21742125
$ws.Range("L77").Value = 45797.31
$ws.Range("M77").Value = -21737757
$ws.Range("N77").Value = -54533.31
$ws.Range("H116").Value = 2674228.2
$ws.Range("I116").Value = 491.5
$ws.Range("J116").Value = 5882712.5
$ws.Range("K116").Value = 491.5
$ws.Range("L116").Value = 5882712.5
$ws.Range("M116").Value = 1802.5
$ws.Range("N116").Value = -5887300.5
$ws.Range("H122").Value = 1410.5625
$ws.Range("I122").Value = 1136.75
$ws.Range("J122").Value = 2232
$ws.Range("K122").Value = 3410.25
$ws.Range("L122").Value = 6696
$ws.Range("M122").Value = -960.25
$ws.Range("N122").Value = -11596
$ws.Range("H123").Value = 54980
$ws.Range("J123").Value = 54980
$ws.Range("L123").Value = 54980
$ws.Range("N123").Value = -64780
$ws.Range("H132").Value = 986687
$ws.Range("I132").Value = 1113786.6
$ws.Range("J132").Value = 1665
$ws.Range("K132").Value = 3341359.8
$ws.Range("L132").Value = 4995
$ws.Range("M132").Value = -3338829.8
$ws.Range("N132").Value = -10055
$ws.Range("H136").Value = 1913.4
$ws.Range("I136").Value = 1682.8889
$ws.Range("J136").Value = 2506.1428
$ws.Range("K136").Value = 5048.6667
$ws.Range("L136").Value = 7518.428400000001
$ws.Range("M136").Value = -2498.6667
$ws.Range("N136").Value = -12618.4284
$ws.Range("H139").Value = 55200.5
$ws.Range("J139").Value = 55200.5
$ws.Range("L139").Value = 55200.5
$ws.Range("N139").Value = -65480.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2674228.2
$ws.Range("I3").Value = 491.5
$ws.Range("J3").Value = 5882712.5
$ws.Range("K3").Value = 491.5
$ws.Range("L3").Value = 5882712.5
$ws.Range("M3").Value = -377.5
$ws.Range("N3").Value = -5882940.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1524.4584
$ws.Range("I31").Value = 1040.7059
$ws.Range("J31").Value = 2699.2856
$ws.Range("K31").Value = 1040.7059
$ws.Range("L31").Value = 2699.2856
$ws.Range("M31").Value = -745.7058999999999
$ws.Range("N31").Value = -3289.2856
$ws.Range("H34").Value = 1524.4584
$ws.Range("I34").Value = 1040.7059
$ws.Range("J34").Value = 2699.2856
$ws.Range("K34").Value = 1040.7059
$ws.Range("L34").Value = 2699.2856
$ws.Range("M34").Value = -838.7058999999999
$ws.Range("N34").Value = -3103.2856
$ws.Range("H58").Value = 2744.3125
$ws.Range("I58").Value = 1566.1428
$ws.Range("J58").Value = 3660.6667
$ws.Range("K58").Value = 1566.1428
$ws.Range("L58").Value = 3660.6667
$ws.Range("M58").Value = -1363.1428
$ws.Range("N58").Value = -4066.6667
$ws.Range("H92").Value = 34996.668
$ws.Range("J92").Value = 34996.668
$ws.Range("L92").Value = 34996.668
$ws.Range("N92").Value = -39988.668
$ws.Range("H136").Value = 2744.3125
$ws.Range("I136").Value = 1566.1428
$ws.Range("J136").Value = 3660.6667
$ws.Range("K136").Value = 4698.428400000001
$ws.Range("L136").Value = 10982.0001
$ws.Range("M136").Value = -2148.428400000001
$ws.Range("N136").Value = -16082.0001
$ws.Range("H141").Value = 19699.8
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 19699.8
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 19699.8
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = -30059.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 5288
$ws.Range("I102").Value = 1005.3333
$ws.Range("K102").Value = 3015.9999
$ws.Range("M102").Value = -581.9998999999998
$ws.Range("H107").Value = 44251.78
$ws.Range("I107").Value = 28400.916
$ws.Range("J107").Value = 101314.9
$ws.Range("K107").Value = 85202.74800000001
$ws.Range("L107").Value = 303944.7
$ws.Range("M107").Value = -83282.74800000001
$ws.Range("N107").Value = -307784.7
$ws.Range("H113").Value = 422.8393
$ws.Range("I113").Value = 408.33334
$ws.Range("J113").Value = 456.11765
$ws.Range("K113").Value = 1225.00002
$ws.Range("L113").Value = 1368.35295
$ws.Range("M113").Value = 944.9999800000001
$ws.Range("N113").Value = -5708.35295
$ws.Range("H131").Value = 1371637.5
$ws.Range("I131").Value = 1035.6428
$ws.Range("J131").Value = 1696865
$ws.Range("K131").Value = 3106.9284
$ws.Range("L131").Value = 5090595
$ws.Range("M131").Value = 1933.0716
$ws.Range("N131").Value = -5100675
$ws.Range("H141").Value = 71432130
$ws.Range("I141").Value = 71432130
$ws.Range("K141").Value = 214296390
$ws.Range("M141").Value = -214291210
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1106.8928
$ws.Range("I97").Value = 1111.32
$ws.Range("J97").Value = 1070
$ws.Range("K97").Value = 1111.32
$ws.Range("L97").Value = 1070
$ws.Range("M97").Value = -615.3199999999999
$ws.Range("N97").Value = -2062
$ws.Range("H102").Value = 1543
$ws.Range("I102").Value = 1660.8889
$ws.Range("J102").Value = 1436.9
$ws.Range("K102").Value = 1660.8889
$ws.Range("L102").Value = 1436.9
$ws.Range("M102").Value = -38.88889999999992
$ws.Range("N102").Value = -4680.9
$ws.Range("H122").Value = 3838.5667
$ws.Range("I122").Value = 4871.1577
$ws.Range("J122").Value = 2055
$ws.Range("K122").Value = 14613.4731
$ws.Range("L122").Value = 6165
$ws.Range("M122").Value = -12163.4731
$ws.Range("N122").Value = -11065
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6795.108
$ws.Range("I132").Value = 8925.950000000001
$ws.Range("J132").Value = 4288.2354
$ws.Range("K132").Value = 26777.85
$ws.Range("L132").Value = 12864.7062
$ws.Range("M132").Value = -24247.85
$ws.Range("N132").Value = -17924.7062
$ws.Range("H137").Value = 61215
$ws.Range("J137").Value = 61215
$ws.Range("L137").Value = 61215
$ws.Range("N137").Value = -71415
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 11123433
$ws.Range("I96").Value = 25001050
$ws.Range("J96").Value = 21340
$ws.Range("K96").Value = 25001050
$ws.Range("L96").Value = 21340
$ws.Range("M96").Value = -24999677
$ws.Range("N96").Value = -24086
$ws.Range("H133").Value = 40447.5
$ws.Range("J133").Value = 40447.5
$ws.Range("L133").Value = 40447.5
$ws.Range("N133").Value = -50567.5
$ws.Range("H136").Value = 68059.28
$ws.Range("I136").Value = 14604.533
$ws.Range("J136").Value = 335333
$ws.Range("K136").Value = 43813.599
$ws.Range("L136").Value = 1005999
$ws.Range("M136").Value = -41263.599
$ws.Range("N136").Value = -1011099
